$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextCell 2 4 '28.436.70'
Set-TextCell 2 5 '  -3.54%  '

Set-TextCell 3 4 '1.956.48'
Set-TextCell 3 5 '  -1.83%  '

Set-TextCell 4 4 '1.007'
Set-TextCell 4 5 '  -0.78%  '

Set-TextCell 5 4 '321.27'
Set-TextCell 5 5 '  -2.36%  '

Set-TextCell 6 5 '  -0.61%  '

Set-TextCell 7 4 '0.4759'
Set-TextCell 7 5 '  -5.13%  '

Set-TextCell 8 4 '0.4052'
Set-TextCell 8 5 '  -4.18%  '

Set-TextCell 9 5 '  -1.06%  '

Set-TextCell 10 4 '0.08406'
Set-TextCell 10 5 '  -5.85%  '

Set-TextCell 11 5 '  -4.70%  '

Set-TextCell 12 4 '22.30'
Set-TextCell 12 5 '  -3.79%  '

Set-TextCell 13 4 '1.948.26'
Set-TextCell 13 5 '  -2.98%  '

Set-TextCell 14 4 '7.609'

Set-TextCell 15 4 '6.147'
Set-TextCell 15 5 '  -4.75%  '

Set-TextCell 16 4 '1.009'
Set-TextCell 16 5 '  -0.62%  '

Set-TextCell 17 4 '90.11'
Set-TextCell 17 5 '  -4.40%  '

Set-TextCell 18 4 '0.00001067'
Set-TextCell 18 5 '  -4.02%  '

Set-TextCell 19 4 '0.06592'
Set-TextCell 19 5 '  -2.27%  '

Set-TextCell 20 4 '18.51'
Set-TextCell 20 5 '  -4.48%  '

Set-TextCell 21 5 '  -0.63%  '

Set-TextCell 22 4 '5.821'
Set-TextCell 22 5 '  -1.77%  '

Set-TextCell 23 4 '28.451.93'
Set-TextCell 23 5 '  -3.60%  '

Set-TextCell 24 4 '11.53'
Set-TextCell 24 5 '  -4.65%  '

Set-TextCell 25 4 '2.290'
Set-TextCell 25 5 '  -1.32%  '

Set-TextCell 26 4 '2.176.35'
Set-TextCell 26 5 '  -3.10%  '

Set-TextCell 27 5 '  -1.16%  '

Set-TextCell 28 4 '20.18'
Set-TextCell 28 5 '  -2.94%  '

Set-TextCell 29 4 '5.916'
Set-TextCell 29 5 '  -6.06%  '

Set-TextCell 30 4 '2.153'
Set-TextCell 30 5 '  -6.40%  '

Set-TextCell 31 4 '123.46'
Set-TextCell 31 5 '  -3.31%  '

Set-TextCell 32 4 '0.9781'
Set-TextCell 32 5 '  -7.55%  '

Set-TextCell 33 4 '0.09593'
Set-TextCell 33 5 '  -3.37%  '

Set-TextCell 34 4 '1.448'
Set-TextCell 34 5 '  -6.38%  '

Set-TextCell 35 4 '5.599'
Set-TextCell 35 5 '  -3.89%  '

Set-TextCell 36 4 '3.658'
Set-TextCell 36 5 '  -3.52%  '

Set-TextCell 37 2 'FraxShare'
Set-TextCell 37 3 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextCell 37 4 '8.905'
Set-TextCell 37 5 '  -3.31%  '

Set-TextCell 38 2 'VeChain'
Set-TextCell 38 3 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell 38 4 '0.02332'
Set-TextCell 38 5 '  -5.16%  '

Set-TextCell 39 4 '0.06209'
Set-TextCell 39 5 '  -2.78%  '

Set-TextCell 40 4 '1.243'
Set-TextCell 40 5 '  -3.96%  '

Set-TextCell 41 4 '0.6202'
Set-TextCell 41 5 '  -4.98%  '

Set-TextCell 42 4 '11.13'
Set-TextCell 42 5 '  -4.00%  '

Set-TextCell 43 5 '  -0.62%  '

Set-TextCell 44 4 '0.1921'
Set-TextCell 44 5 '  -5.72%  '

Set-TextCell 45 4 '1.358'
Set-TextCell 45 5 '  +3.72%  '

Set-TextCell 46 4 '0.5952'
Set-TextCell 46 5 '  -5.92%  '

Set-TextCell 47 4 '13.04'
Set-TextCell 47 5 '  -2.97%  '

Set-TextCell 48 4 '2.058'
Set-TextCell 48 5 '  -6.82%  '

Set-TextCell 49 4 '3.393'
Set-TextCell 49 5 '  -3.09%  '

Set-TextCell 50 5 '  -3.16%  '

Set-TextCell 51 4 '0.06827'
Set-TextCell 51 5 '  -1.81%  '

